$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Step 1: bump the date in A1 by one day (45308 -> 45309)
$ws.Range("A1").Value = 45309

# Step 2: update the price list values in column D
$ws.Range("D29").Value = 348.194
$ws.Range("D30").Value = 368.347
$ws.Range("D31").Value = 396.64
